# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. This updates the DAMSLTag (column I) and
# DialogAct (column J) values for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> (DAMSLTag, DialogAct)
$updates = @{
    9   = @('sd', 'Statement-non-opinion')
    53  = @('ba', 'Appreciation')
    58  = @('ba', 'Appreciation')
    75  = @('%', 'Uninterpretable')
    87  = @('sd', 'Statement-non-opinion')
    94  = @('b', 'Acknowledge (Backchannel)')
    108 = @('b', 'Acknowledge (Backchannel)')
    124 = @('%', 'Uninterpretable')
    125 = @('%', 'Uninterpretable')
    146 = @('aa', 'Agree/Accept')
    164 = @('sd', 'Statement-non-opinion')
    180 = @('sd', 'Statement-non-opinion')
    184 = @('%', 'Uninterpretable')
    219 = @('b', 'Acknowledge (Backchannel)')
    227 = @('sv', 'Statement-opinion')
    228 = @('sd', 'Statement-non-opinion')
    230 = @('sd', 'Statement-non-opinion')
    246 = @('sv', 'Statement-opinion')
    263 = @('sd', 'Statement-non-opinion')
    281 = @('sv', 'Statement-opinion')
    283 = @('sd', 'Statement-non-opinion')
    295 = @('%', 'Uninterpretable')
    297 = @('aa', 'Agree/Accept')
    305 = @('ba', 'Appreciation')
    308 = @('b', 'Acknowledge (Backchannel)')
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    $ws.Cells.Item($row, 9).Value = $values[0]
    $ws.Cells.Item($row, 10).Value = $values[1]
}
